$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) columns for each coin row
# per the refreshed cryptos list feed.

$ws.Range("D2").Value = "29.507.69"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").Value = "1.917.03"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.45"
$ws.Range("E5").Value = "  -1.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  +3.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4095"
$ws.Range("E8").Value = "  -0.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.95"
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08039"
$ws.Range("E10").Value = "  +0.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.010"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.41"
$ws.Range("E12").Value = "  +2.76%  "

$ws.Range("D13").Value = "1.908.12"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.942"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.149"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.57"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06603"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001031"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "29.512.24"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.547"
$ws.Range("E23").Value = "  +2.09%  "

$ws.Range("E24").Value = "  +1.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.204"
$ws.Range("E25").Value = "  -1.26%  "

$ws.Range("D26").Value = "2.127.12"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.98"
$ws.Range("E27").Value = "  -2.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.83"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.822"
$ws.Range("E29").Value = "  +7.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.134"
$ws.Range("E30").Value = "  +0.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.72"
$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.053"
$ws.Range("E32").Value = "  +7.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09540"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("E34").Value = "  -0.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.573"
$ws.Range("E35").Value = "  -0.77%  "

$ws.Range("E36").Value = "  +1.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06108"
$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02252"
$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.324"
$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.174"
$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5879"
$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("E42").Value = "  +9.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1844"
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.13"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.08029"
$ws.Range("E45").Value = "  +13.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.286"
$ws.Range("E46").Value = "  +1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.929"
$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.19"
$ws.Range("E50").Value = "  +1.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.54"
$ws.Range("E51").Value = "  -4.89%  "

# Rows 47/48: EnergySwap and Decentraland swapped positions in the ranking
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5538"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.13"
$ws.Range("E48").Value = "  +0.83%  "
